$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(460, 1).Value = 44830
$ws.Cells.Item(460, 2).Value = 'KA03MR2992'
$ws.Cells.Item(460, 3).Value = 'ALTO K10'
$ws.Cells.Item(460, 4).Value = 'PMS & SUSPENSION'
$ws.Cells.Item(460, 5).Value = 'WORK DONE DELIVERED'
$ws.Cells.Item(460, 6).Value = 13190
$ws.Cells.Item(460, 7).Value = 'CARD'

$ws.Cells.Item(461, 1).Value = 44830
$ws.Cells.Item(461, 2).Value = 'KA13V4655'
$ws.Cells.Item(461, 3).Value = 'SWIFT'
$ws.Cells.Item(461, 4).Value = 'PMS'
$ws.Cells.Item(461, 5).Value = 'WORK DONE DELIVERED'
$ws.Cells.Item(461, 6).Value = 7948
$ws.Cells.Item(461, 7).Value = 'CREDIT'

$ws.Cells.Item(462, 1).Value = 44830
$ws.Cells.Item(462, 2).Value = 'KA03MT0636'
$ws.Cells.Item(462, 3).Value = 'ECOSPORT'
$ws.Cells.Item(462, 4).Value = 'GENERAL CHECKUP'
$ws.Cells.Item(462, 5).Value = 'WORK DONE DELIVERED'
$ws.Cells.Item(462, 6).Value = 1374
$ws.Cells.Item(462, 7).Value = 'CREDIT'

$ws.Cells.Item(463, 1).Value = 44830
$ws.Cells.Item(463, 2).Value = 'KA04MN7139'
$ws.Cells.Item(463, 3).Value = 'I10'
$ws.Cells.Item(463, 4).Value = 'PMS'
$ws.Cells.Item(463, 5).Value = 'WORK DONE DELIVERED'
$ws.Cells.Item(463, 6).Value = 3403
$ws.Cells.Item(463, 7).Value = 'GPAY'

$ws.Cells.Item(464, 1).Value = 44830
$ws.Cells.Item(464, 2).Value = 'KA03NA5054'
$ws.Cells.Item(464, 3).Value = 'POLO'
$ws.Cells.Item(464, 4).Value = 'BODY SHOP'
$ws.Cells.Item(464, 5).Value = 'WORK DONE'
$ws.Cells.Item(464, 6).Value = 40965
$ws.Cells.Item(464, 7).Value = '  INSURANCE'

$ws.Cells.Item(465, 1).Value = 44831
$ws.Cells.Item(465, 2).Value = 'KA04MM4818'
$ws.Cells.Item(465, 3).Value = 'RITZ'
$ws.Cells.Item(465, 4).Value = 'PART SEALS                     WW'
$ws.Cells.Item(465, 5).Value = 'WORK DONE DELIVERED'
$ws.Cells.Item(465, 6).Value = 1136
$ws.Cells.Item(465, 7).Value = 'CREDIT'

$ws.Cells.Item(466, 1).Value = 44831
$ws.Cells.Item(466, 2).Value = 'KA51MM2838'
$ws.Cells.Item(466, 3).Value = 'SPARK'
$ws.Cells.Item(466, 4).Value = 'PMS'
$ws.Cells.Item(466, 5).Value = 'WORK DONE DELIVERED'
$ws.Cells.Item(466, 6).Value = 3618
$ws.Cells.Item(466, 7).Value = 'P PAY'

$ws.Cells.Item(467, 1).Value = 44831
$ws.Cells.Item(467, 2).Value = 'KA53MD8318'
$ws.Cells.Item(467, 3).Value = 'JAZZ'
$ws.Cells.Item(467, 4).Value = 'PMS '
$ws.Cells.Item(467, 5).Value = 'WORK DONE DELIVERED'
$ws.Cells.Item(467, 6).Value = 6380
$ws.Cells.Item(467, 7).Value = 'P PAY'

$ws.Cells.Item(468, 1).Value = 44831
$ws.Cells.Item(468, 2).Value = 'KA53MB1800'
$ws.Cells.Item(468, 3).Value = 'SCALA'
$ws.Cells.Item(468, 4).Value = 'GENERAL CHECKUP '
$ws.Cells.Item(468, 5).Value = 'WORK DONE DELIVERED'
$ws.Cells.Item(468, 6).Value = 708
$ws.Cells.Item(468, 7).Value = 'GPAY'

$ws.Cells.Item(469, 1).Value = 44831
$ws.Cells.Item(469, 2).Value = 'KA01MJ3412'
$ws.Cells.Item(469, 3).Value = 'VISTA'
$ws.Cells.Item(469, 4).Value = 'CLUTCH PROBLEM'
$ws.Cells.Item(469, 5).Value = 'WORK DONE DELIVERED'
$ws.Cells.Item(469, 6).Value = 13314
$ws.Cells.Item(469, 7).Value = 'P PAY'

$ws.Cells.Item(470, 1).Value = 44831
$ws.Cells.Item(470, 2).Value = 'AP09BX8688'
$ws.Cells.Item(470, 3).Value = 'RITZ'
$ws.Cells.Item(470, 4).Value = 'PMS'
$ws.Cells.Item(470, 5).Value = 'WORK DONE DELIVERED'
$ws.Cells.Item(470, 6).Value = 11912
$ws.Cells.Item(470, 7).Value = 'CREDIT'

$ws.Cells.Item(471, 1).Value = 44831
$ws.Cells.Item(471, 2).Value = 'KA03MN9673'
$ws.Cells.Item(471, 3).Value = 'POLO'
$ws.Cells.Item(471, 4).Value = 'RUNNING REPAIR'
$ws.Cells.Item(471, 5).Value = 'WORK IN PROGRESS'

$ws.Cells.Item(472, 1).Value = 44831
$ws.Cells.Item(472, 2).Value = 'KA51MK9302'
$ws.Cells.Item(472, 3).Value = 'CRETA'
$ws.Cells.Item(472, 4).Value = 'PMS'
$ws.Cells.Item(472, 5).Value = 'WORK DONE '

$ws.Cells.Item(473, 1).Value = 44831
$ws.Cells.Item(473, 2).Value = 'KA03MK9302'
$ws.Cells.Item(473, 3).Value = 'SPARK'
$ws.Cells.Item(473, 4).Value = 'PMS'
$ws.Cells.Item(473, 5).Value = 'WORK IN PROGRESS'

# Update the active selection/view to match the post-edit state
$ws.Activate()
$ws.Range("A474").Select()
$excel.ActiveWindow.ScrollRow = 451
$excel.ActiveWindow.ScrollColumn = 1
